$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple numeric value updates (style unchanged) ---
$ws.Range("F14").Value = 1
$ws.Range("L15").Value = -42.105263157894
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = -37.5
$ws.Range("L16").Value = 25.892857142857
$ws.Range("M16").Value = -43.373493975903
$ws.Range("N16").Value = -90.808344198174
$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 125
$ws.Range("F17").Value = 25
$ws.Range("G17").Value = 30
$ws.Range("H17").Value = -16.666666666666
$ws.Range("I17").Value = 317
$ws.Range("J17").Value = 283
$ws.Range("K17").Value = 12.014134275618
$ws.Range("L17").Value = 7.094594594594
$ws.Range("M17").Value = 14.440433212996
$ws.Range("N17").Value = -62.749706227967
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 7
$ws.Range("H18").Value = -41.666666666666
$ws.Range("I18").Value = 122
$ws.Range("J18").Value = 88
$ws.Range("K18").Value = 38.636363636363
$ws.Range("L18").Value = -2.4
$ws.Range("M18").Value = -57.042253521126
$ws.Range("N18").Value = -92.848769050410
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 75
$ws.Range("F19").Value = 37
$ws.Range("G19").Value = 24
$ws.Range("H19").Value = 54.166666666666
$ws.Range("I19").Value = 328
$ws.Range("J19").Value = 246
$ws.Range("K19").Value = 33.333333333333
$ws.Range("L19").Value = 28.627450980392
$ws.Range("M19").Value = -16.112531969309
$ws.Range("N19").Value = -48.019017432646
$ws.Range("C20").Value = 4
$ws.Range("E20").Value = 300
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = 150
$ws.Range("I20").Value = 131
$ws.Range("J20").Value = 107
$ws.Range("K20").Value = 22.429906542056
$ws.Range("L20").Value = 22.429906542056
$ws.Range("M20").Value = -7.092198581560
$ws.Range("N20").Value = -89.477911646586
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 10
$ws.Range("E21").Value = 120
$ws.Range("F21").Value = 85
$ws.Range("G21").Value = 80
$ws.Range("H21").Value = 6.25
$ws.Range("I21").Value = 1058
$ws.Range("J21").Value = 850
$ws.Range("K21").Value = 24.470588235294
$ws.Range("L21").Value = 14.875135722041
$ws.Range("M21").Value = -22.263041880969
$ws.Range("N21").Value = -82.561397725399
$ws.Range("M23").Value = -80
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = -4.166666666666
$ws.Range("F24").Value = 98
$ws.Range("G24").Value = 98
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 1009
$ws.Range("J24").Value = 987
$ws.Range("K24").Value = 2.228976697061
$ws.Range("L24").Value = 9.199134199134
$ws.Range("M24").Value = 9.912854030501
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = 100
$ws.Range("F25").Value = 53
$ws.Range("G25").Value = 46
$ws.Range("H25").Value = 15.217391304347
$ws.Range("I25").Value = 521
$ws.Range("J25").Value = 430
$ws.Range("K25").Value = 21.162790697674
$ws.Range("L25").Value = 50.578034682080
$ws.Range("M25").Value = -31.357048748353
$ws.Range("L26").Value = -24.137931034482
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 9
$ws.Range("H27").Value = -55.555555555555
$ws.Range("J27").Value = 62
$ws.Range("K27").Value = -12.903225806451
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = -50
$ws.Range("M28").Value = -42.105263157894
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = -50
$ws.Range("M29").Value = -33.333333333333

# --- Numeric -> text "0" (shared string, style 14), source C14 ---
$ws.Range("C14").Copy($ws.Range("C16"))
$ws.Range("C14").Copy($ws.Range("D16"))
$ws.Range("C14").Copy($ws.Range("D22"))
$ws.Range("C14").Copy($ws.Range("D26"))
$ws.Range("C14").Copy($ws.Range("C27"))
$ws.Range("C14").Copy($ws.Range("D28"))
$ws.Range("C14").Copy($ws.Range("D29"))

# --- Numeric -> text "***.*" (shared string, style 14), source E14 ---
$ws.Range("E14").Copy($ws.Range("E16"))
$ws.Range("E14").Copy($ws.Range("E22"))
$ws.Range("E14").Copy($ws.Range("E26"))
$ws.Range("E14").Copy($ws.Range("E28"))
$ws.Range("E14").Copy($ws.Range("E29"))

# --- Text -> numeric (copy format only, then set value) ---
$ws.Range("F16").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").Value = 3
$ws.Range("H16").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E27").Value = -100